$d = $word.ActiveDocument

# Update the date/day heading at the top of the document
$d.Content.Find.Execute("2024-05-16 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-05-17 Friday", 2) | Out-Null

# New values for each of the 100 arithmetic-problem cells, in row-major
# (left-to-right, top-to-bottom) order matching the table layout.
$newValues = @(
    "3+49=",
    "4+88=",
    "86-58=",
    "75+6=",
    "45+17=",
    "30-29=",
    "6+26=",
    "37+9=",
    "7+84=",
    "47+29=",
    "39+42=",
    "34-16=",
    "52-18=",
    "48+15=",
    "70-59=",
    "5+9=",
    "77+16=",
    "5+16=",
    "9+12=",
    "3+8=",
    "50-39=",
    "93-65=",
    "57+4=",
    "8+5=",
    "84-79=",
    "33+38=",
    "70-61=",
    "63-29=",
    "15+36=",
    "77+4=",
    "57+29=",
    "9+55=",
    "60-33=",
    "39+37=",
    "45+8=",
    "81-29=",
    "44+48=",
    "80-63=",
    "94-15=",
    "80-54=",
    "25-16=",
    "7+54=",
    "28+26=",
    "80-49=",
    "83-4=",
    "7+79=",
    "15+26=",
    "85+6=",
    "25+68=",
    "19+29=",
    "46+48=",
    "63-34=",
    "85-38=",
    "25+67=",
    "52-43=",
    "91-32=",
    "16+19=",
    "60-52=",
    "94-36=",
    "3+28=",
    "66-59=",
    "25+17=",
    "58+23=",
    "6+46=",
    "57+26=",
    "32-25=",
    "91-23=",
    "81-12=",
    "31-24=",
    "14+39=",
    "86-27=",
    "40-14=",
    "93-86=",
    "92-17=",
    "15+36=",
    "63-38=",
    "3+78=",
    "72-16=",
    "90-7=",
    "43-14=",
    "50-34=",
    "39+18=",
    "28+33=",
    "90-38=",
    "83-58=",
    "54+19=",
    "43+8=",
    "58+38=",
    "37+28=",
    "77+8=",
    "91-25=",
    "13+48=",
    "61-43=",
    "39+38=",
    "90-34=",
    "21-15=",
    "7+38=",
    "62-9=",
    "81-26=",
    "54-37="
)

$t = $d.Tables.Item(1)
$rows = $t.Rows.Count
$cols = $t.Columns.Count

$idx = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $t.Cell($r, $c)
        $cellRange = $cell.Range
        # Trim the trailing paragraph mark + cell mark (2 chars) off the cell range
        $cellRange.MoveEnd(1, -2) | Out-Null
        $cellRange.Text = $newValues[$idx]
        $idx = $idx + 1
    }
}

Write-Host "Done: updated $idx cells"
